$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates: extraction re-run landed on a different requisition/patient ---

# Requisition number
$ws.Range("A2").Value = "REQUISITION:P2418L001Y"

# City (format as text so it's not mangled)
$ws.Range("E2").Value = "lowa Park, TX 763678633"

# ZipCode - force text so the all-digit string isn't coerced to a number
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "763678633"
$ws.Range("F2").Style = "Normal"

# HomePhone
$ws.Range("G2").Value = "940-704-9644"

# Phys Address
$ws.Range("H2").Value = "4327 Barnett Road Wichita Falls. TX 763102303"

# Patient Name
$ws.Range("I2").Value = "Redclift, Reynold"

# Gender
$ws.Range("J2").Value = "L (M/78)"

# DateOfBirth - force text so the date-like string isn't coerced to a serial date
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "1945-05-23"
$ws.Range("K2").Style = "Normal"

# Address
$ws.Range("L2").Value = "11681 Longley Road"

# MRN / External MRN
$ws.Range("M2").Value = "RERE0001"
$ws.Range("N2").Value = "RERE0001"

# Primary Ins Num
$ws.Range("P2").Value = "MR005) (2W51VF0GN34)"

# SecondaryInsurance_CompanyName / Num
$ws.Range("Q2").Value = "Mutual of Omaha Insurance Company"
$ws.Range("R2").Value = "(72861188)"

# Site Location 1 / Physical Exam 1
$ws.Range("T2").Value = "Right Suprascapular Back = Lateral"
$ws.Range("W2").Value = "Pink papule (Right Suprascapular Back - Lateral)"

# Site Location 2 / Other 2 / Type Of Procedure 2 / Physical Exam 2
# (previously this second page of the dialogue was only partly filled in)
$ws.Range("X2").Value = "Right Ear = Descending Helix"
$ws.Range("Y2").Value = "Helix Neoplasm of uncertain behavior of skin"
$ws.Range("Z2").Value = "Biopsy (Tangential (Shave))"
$ws.Range("AA2").Value = "Exam: Pink papule (Right Ear = Descending Helix)"
